# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" (cloned from "2021-Q4" so header text,
#   column layout and cell styling match) positioned right before "总计",
#   then overwrite its data rows with the 2022-Q1 fund-holding figures.
# - Update the "总计" (totals) summary sheet: add a new top data row for
#   "2022-Q1" and push the existing "2021-Q4"/"2021-Q3" rows down by one,
#   re-numbering the index column (A) accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell as genuine text (no quote-prefix /
# number auto-coercion residue left behind on the cell's style).
# ---------------------------------------------------------------------
function Set-TextCell($ws, $addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet by cloning "2021-Q4" (same header row /
#    column widths / cell styles), positioned immediately before "总计".
# ---------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore = $wb.Worksheets.Item("总计")
$q4Sheet.Copy($totalSheetBefore)

# Re-resolve "总计" fresh (the pre-copy reference's cached .Index goes
# stale once the sheet collection is mutated) so we can reliably find the
# sheet that was just inserted directly ahead of it.
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item($totalSheet.Index - 1)
$newSheet.Name = "2022-Q1"

# Overwrite the cloned data rows with the actual 2022-Q1 figures.
Set-TextCell $newSheet "B2" "005262"
Set-TextCell $newSheet "C2" "鑫元欣享灵活配置混合A"
Set-TextCell $newSheet "D2" "0.89"
Set-TextCell $newSheet "E2" "85.87"
Set-TextCell $newSheet "F2" "2.89"
Set-TextCell $newSheet "G2" "0.0257"
$newSheet.Range("H2").Value = 10

Set-TextCell $newSheet "B3" "005263"
Set-TextCell $newSheet "C3" "鑫元欣享灵活配置混合C"
Set-TextCell $newSheet "D3" "0.09"
Set-TextCell $newSheet "E3" "85.87"
Set-TextCell $newSheet "F3" "2.89"
Set-TextCell $newSheet "G3" "0.0026"
$newSheet.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2. Update "总计": push the two existing rows down one slot (copying
#    cells so the original per-row styling/number formats move with the
#    data) and write the new "2022-Q1" row on top.
# ---------------------------------------------------------------------
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A4:D4"))
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.03

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# ---------------------------------------------------------------------
# 3. Restore the originally-active tab ("2021-Q3"): copying a sheet
#    activates the copy, so without this the saved file would mark the
#    new "2022-Q1" tab as selected instead of leaving "2021-Q3" as-is.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()

